# Append the new run-log row (row 24) that records the 2025-08-17 09:37:58 UTC
# "SKIPPED" run, mirroring the formatting of the previous row (row 23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 23
$newRow = 24

# Copy the formatting (cell style) of the last existing row onto the new row
# so every new cell (including the blank ones) ends up styled like the rest
# of the log (centered alignment, s="3").
$ws.Range("A" + $lastRow + ":H" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":H" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new run's data.
$ws.Range("A" + $newRow).Value2 = "2025-08-17 09:37:58 UTC"
$ws.Range("B" + $newRow).Value2 = "2025-08-17 15:07:58 IST"
$ws.Range("C" + $newRow).Value2 = "SKIPPED"
$ws.Range("D" + $newRow).Value2 = "No change in PDF. Skipping download & Excel update."
$ws.Range("E" + $newRow).Value2 = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Range("G" + $newRow).Value2 = 0
